$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new pin entry: Pressure Transducer on pin 3 (row 10)
$ws.Range("A10").Value = "PRESSURE TRANSDUCER"
$ws.Range("B10").Value = 3

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("C18").Select()
